$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("F3").Value = -4
$ws.Range("F8").Value = 3
$ws.Range("F9").Value = 0
$ws.Range("F14").Value = 6
$ws.Range("F15").Value = 0
$ws.Range("F16").Value = 0
$ws.Range("F18").Value = -3
